# Reorder the "<Field>, <Degree(s)>" education lines to "<Degree(s)>, <Field>".
# The target runs keep identical run-level formatting (Times New Roman,
# italic, sz 20) to the original single run, so a plain text Find/Replace
# (which leaves the run's rPr untouched) reproduces the rendered/semantic
# result of the diff.

$d = $word.ActiveDocument

# Ph.D., M.S. line (Georgia Tech row)
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(
    "Mechanical Engineering, Ph.D., M.S.", $true, $false, $false, $false, $false,
    $true, 1, $false, "Ph.D., M.S., Mechanical Engineering", 2)
Write-Host "Replaced PhD/MS line: $found1"

# B.S. line (Pittsburgh row)
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    "Mechanical Engineering, B.S.", $true, $false, $false, $false, $false,
    $true, 1, $false, "B.S., Mechanical Engineering", 2)
Write-Host "Replaced BS line: $found2"
